$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block of renamed-flat entries, appended after the existing table
# (mirrors the "Name" header / path / old-new pairs layout already on the sheet).
$ws.Range("A57").Value = "Name"
$ws.Range("D57").Value = "Name"

$ws.Range("A59").Value = "Z:\BLACK\lexicon\pk3\Textures\renamed\2\"
$ws.Range("D59").Value = "Z:\BLACK\lexicon\pk3\Textures\renamed\2\"

$ws.Range("A60").Value = "CEIL1_2"
$ws.Range("D60").Value = "_T100050"

$ws.Range("A61").Value = "GATE1"
$ws.Range("D61").Value = "_T100051"

$ws.Range("A62").Value = "GATE2"
$ws.Range("D62").Value = "_T100052"

$ws.Range("A63").Value = "GATE3"
$ws.Range("D63").Value = "_T100053"

$ws.Range("A64").Value = "GATE4"
$ws.Range("D64").Value = "_T100054"

$ws.Range("A65").Value = "RROCK01"
$ws.Range("D65").Value = "_T100055"

$ws.Range("A66").Value = "RROCK02"
$ws.Range("D66").Value = "_T100056"

# Selection ends on the last populated cell, matching the authored scroll state.
$ws.Range("D66").Select()
